$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-11 with new dates / litros / odm values ---
$ws.Range("A2").Value = 45536
$ws.Range("D2").Value = 43.250999999999998
$ws.Range("E2").Value = 96349.259000000005

$ws.Range("A3").Value = 45537
$ws.Range("D3").Value = 45.222000000000001
$ws.Range("E3").Value = 96756.257000000012

$ws.Range("A4").Value = 45538
$ws.Range("D4").Value = 47.253
$ws.Range("E4").Value = 97181.534000000014

$ws.Range("A5").Value = 45539
$ws.Range("D5").Value = 44.529000000000003
$ws.Range("E5").Value = 97582.295000000013

$ws.Range("A6").Value = 45540
$ws.Range("D6").Value = 42.369
$ws.Range("E6").Value = 97963.616000000009

$ws.Range("A7").Value = 45541
$ws.Range("D7").Value = 41.253999999999998
$ws.Range("E7").Value = 96830.032000000007

$ws.Range("A8").Value = 45542
$ws.Range("D8").Value = 42.057000000000002
$ws.Range("E8").Value = 97166.488000000012

$ws.Range("A9").Value = 45543
$ws.Range("D9").Value = 39.546999999999997
$ws.Range("E9").Value = 97482.864000000016

$ws.Range("A10").Value = 45544
$ws.Range("D10").Value = 37.027000000000001
$ws.Range("E10").Value = 97779.080000000016

$ws.Range("A11").Value = 45545
$ws.Range("D11").Value = 42.658000000000001
$ws.Range("E11").Value = 98120.344000000012

# --- Rows 12-14 no longer hold data; wipe their contents but keep styling ---
$ws.Range("A12:G14").ClearContents()

# --- Move the sheet's active selection to H16 ---
$ws.Range("H16").Select()
